$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.455362044514542
$ws.Cells.Item(2, 3).Value = 1.655778082260271
$ws.Cells.Item(2, 4).Value = 3.537761648806719
$ws.Cells.Item(2, 5).Value = 10.19245300693656
$ws.Cells.Item(2, 7).Value = 16.84135478251809
$ws.Cells.Item(3, 2).Value = 3.286832544864788
$ws.Cells.Item(3, 3).Value = 1.655778082260271
$ws.Cells.Item(3, 4).Value = 0.7527432677738641
$ws.Cells.Item(3, 5).Value = 0.4942365360607697
$ws.Cells.Item(3, 7).Value = 6.189590430959694
$ws.Cells.Item(4, 2).Value = 3.286832544864788
$ws.Cells.Item(4, 3).Value = 1.655778082260271
$ws.Cells.Item(4, 4).Value = 261.3203778131603
$ws.Cells.Item(4, 5).Value = 10.19245300693656
$ws.Cells.Item(4, 7).Value = 276.4554414472219
$ws.Cells.Item(5, 2).Value = 1.455362044514542
$ws.Cells.Item(5, 3).Value = 1.655778082260271
$ws.Cells.Item(5, 4).Value = 0.7527432677738641
$ws.Cells.Item(5, 5).Value = 0.4942365360607697
$ws.Cells.Item(5, 7).Value = 4.358119930609447
$ws.Cells.Item(6, 2).Value = 3.286832544864788
$ws.Cells.Item(6, 3).Value = 1.655778082260271
$ws.Cells.Item(6, 4).Value = 22.3905356188092
$ws.Cells.Item(6, 5).Value = 0.4942365360607697
$ws.Cells.Item(6, 7).Value = 27.82738278199502
$ws.Cells.Item(7, 2).Value = 1.455362044514542
$ws.Cells.Item(7, 3).Value = 1.655778082260271
$ws.Cells.Item(7, 4).Value = 3.537761648806719
$ws.Cells.Item(7, 5).Value = 0.4942365360607697
$ws.Cells.Item(7, 7).Value = 7.143138311642302
$ws.Cells.Item(8, 2).Value = 0.1190320826869504
$ws.Cells.Item(8, 3).Value = 0.306821227259698
$ws.Cells.Item(8, 4).Value = 0.7527432677738641
$ws.Cells.Item(8, 5).Value = 0.4942365360607697
$ws.Cells.Item(8, 7).Value = 1.672833113781282
$ws.Cells.Item(9, 2).Value = 0.2917716402565462
$ws.Cells.Item(9, 3).Value = 0.306821227259698
$ws.Cells.Item(9, 4).Value = 0.1494219747398047
$ws.Cells.Item(9, 5).Value = 0.4942365360607697
$ws.Cells.Item(9, 7).Value = 1.242251378316819
$ws.Cells.Item(10, 2).Value = 0.6606524410359556
$ws.Cells.Item(10, 3).Value = 1.655778082260271
$ws.Cells.Item(10, 4).Value = 0.7527432677738641
$ws.Cells.Item(10, 5).Value = 0.4942365360607697
$ws.Cells.Item(10, 7).Value = 3.56341032713086
$ws.Cells.Item(11, 2).Value = 3.286832544864788
$ws.Cells.Item(11, 3).Value = 1.655778082260271
$ws.Cells.Item(11, 4).Value = 0.7527432677738641
$ws.Cells.Item(11, 5).Value = 0.4942365360607697
$ws.Cells.Item(11, 7).Value = 6.189590430959694
$ws.Cells.Item(12, 2).Value = 3.286832544864788
$ws.Cells.Item(12, 3).Value = 1.655778082260271
$ws.Cells.Item(12, 4).Value = 3.537761648806719
$ws.Cells.Item(12, 5).Value = 0.4942365360607697
$ws.Cells.Item(12, 7).Value = 8.974608811992548
$ws.Cells.Item(13, 2).Value = 3.286832544864788
$ws.Cells.Item(13, 3).Value = 1.655778082260271
$ws.Cells.Item(13, 4).Value = 0.7527432677738641
$ws.Cells.Item(13, 5).Value = 0.4942365360607697
$ws.Cells.Item(13, 7).Value = 6.189590430959694
$ws.Cells.Item(14, 2).Value = 0.04271373187048222
$ws.Cells.Item(14, 3).Value = 0.306821227259698
$ws.Cells.Item(14, 4).Value = 0.7527432677738641
$ws.Cells.Item(14, 5).Value = 0.4942365360607697
$ws.Cells.Item(14, 7).Value = 1.596514762964814
$ws.Cells.Item(15, 2).Value = 0.6606524410359556
$ws.Cells.Item(15, 3).Value = 1.655778082260271
$ws.Cells.Item(15, 4).Value = 22.3905356188092
$ws.Cells.Item(15, 5).Value = 10.19245300693656
$ws.Cells.Item(15, 7).Value = 34.89941914904198
$ws.Cells.Item(16, 2).Value = 0.1190320826869504
$ws.Cells.Item(16, 3).Value = 117.745847958593
$ws.Cells.Item(16, 4).Value = 0.7527432677738641
$ws.Cells.Item(16, 5).Value = 10.19245300693656
$ws.Cells.Item(16, 7).Value = 128.8100763159904
$ws.Cells.Item(17, 2).Value = 3.286832544864788
$ws.Cells.Item(17, 3).Value = 1.655778082260271
$ws.Cells.Item(17, 4).Value = 0.1494219747398047
$ws.Cells.Item(17, 5).Value = 0.4942365360607697
$ws.Cells.Item(17, 7).Value = 5.586269137925634
$ws.Cells.Item(18, 2).Value = 1.455362044514542
$ws.Cells.Item(18, 3).Value = 1.655778082260271
$ws.Cells.Item(18, 4).Value = 3.537761648806719
$ws.Cells.Item(18, 5).Value = 0.4942365360607697
$ws.Cells.Item(18, 7).Value = 7.143138311642302
$ws.Cells.Item(19, 2).Value = 3.286832544864788
$ws.Cells.Item(19, 3).Value = 1.655778082260271
$ws.Cells.Item(19, 4).Value = 3.537761648806719
$ws.Cells.Item(19, 5).Value = 0.4942365360607697
$ws.Cells.Item(19, 7).Value = 8.974608811992548
$ws.Cells.Item(20, 2).Value = 3.286832544864788
$ws.Cells.Item(20, 3).Value = 1.655778082260271
$ws.Cells.Item(20, 4).Value = 3.537761648806719
$ws.Cells.Item(20, 5).Value = 0.4942365360607697
$ws.Cells.Item(20, 7).Value = 8.974608811992548
$ws.Cells.Item(21, 2).Value = 3.286832544864788
$ws.Cells.Item(21, 3).Value = 1.655778082260271
$ws.Cells.Item(21, 4).Value = 22.3905356188092
$ws.Cells.Item(21, 5).Value = 10.19245300693656
$ws.Cells.Item(21, 7).Value = 37.52559925287081
$ws.Cells.Item(22, 2).Value = 1.455362044514542
$ws.Cells.Item(22, 3).Value = 1.655778082260271
$ws.Cells.Item(22, 4).Value = 0.7527432677738641
$ws.Cells.Item(22, 5).Value = 0.4942365360607697
$ws.Cells.Item(22, 7).Value = 4.358119930609447
$ws.Cells.Item(23, 2).Value = 3.286832544864788
$ws.Cells.Item(23, 3).Value = 1.655778082260271
$ws.Cells.Item(23, 4).Value = 0.1494219747398047
$ws.Cells.Item(23, 5).Value = 0.4942365360607697
$ws.Cells.Item(23, 7).Value = 5.586269137925634
$ws.Cells.Item(24, 2).Value = 3.286832544864788
$ws.Cells.Item(24, 3).Value = 10.34677158129881
$ws.Cells.Item(24, 4).Value = 6708.013860684405
$ws.Cells.Item(24, 5).Value = 10.19245300693656
$ws.Cells.Item(24, 7).Value = 6731.839917817505
$ws.Cells.Item(25, 2).Value = 0.04271373187048222
$ws.Cells.Item(25, 3).Value = 0.04071648406533734
$ws.Cells.Item(25, 4).Value = 22.3905356188092
$ws.Cells.Item(25, 5).Value = 0.4942365360607697
$ws.Cells.Item(25, 7).Value = 22.96820237080578
$ws.Cells.Item(26, 2).Value = 3.286832544864788
$ws.Cells.Item(26, 3).Value = 1.655778082260271
$ws.Cells.Item(26, 4).Value = 22.3905356188092
$ws.Cells.Item(26, 5).Value = 0.4942365360607697
$ws.Cells.Item(26, 7).Value = 27.82738278199502
$ws.Cells.Item(27, 2).Value = 1.455362044514542
$ws.Cells.Item(27, 3).Value = 117.745847958593
$ws.Cells.Item(27, 4).Value = 3.537761648806719
$ws.Cells.Item(27, 5).Value = 10.19245300693656
$ws.Cells.Item(27, 7).Value = 132.9314246588508
$ws.Cells.Item(28, 2).Value = 0.6606524410359556
$ws.Cells.Item(28, 3).Value = 1.655778082260271
$ws.Cells.Item(28, 4).Value = 3.537761648806719
$ws.Cells.Item(28, 5).Value = 10.19245300693656
$ws.Cells.Item(28, 7).Value = 16.0466451790395
$ws.Cells.Item(29, 2).Value = 3.286832544864788
$ws.Cells.Item(29, 3).Value = 1.655778082260271
$ws.Cells.Item(29, 4).Value = 0.7527432677738641
$ws.Cells.Item(29, 5).Value = 0.4942365360607697
$ws.Cells.Item(29, 7).Value = 6.189590430959694
$ws.Cells.Item(30, 2).Value = 3.286832544864788
$ws.Cells.Item(30, 3).Value = 1.655778082260271
$ws.Cells.Item(30, 4).Value = 3.537761648806719
$ws.Cells.Item(30, 5).Value = 0.4942365360607697
$ws.Cells.Item(30, 7).Value = 8.974608811992548
$ws.Cells.Item(31, 2).Value = 3.286832544864788
$ws.Cells.Item(31, 3).Value = 1.655778082260271
$ws.Cells.Item(31, 4).Value = 0.1494219747398047
$ws.Cells.Item(31, 5).Value = 0.4942365360607697
$ws.Cells.Item(31, 7).Value = 5.586269137925634
$ws.Cells.Item(32, 2).Value = 3.286832544864788
$ws.Cells.Item(32, 3).Value = 1.655778082260271
$ws.Cells.Item(32, 4).Value = 0.1494219747398047
$ws.Cells.Item(32, 5).Value = 0.4942365360607697
$ws.Cells.Item(32, 7).Value = 5.586269137925634
$ws.Cells.Item(33, 2).Value = 0.6606524410359556
$ws.Cells.Item(33, 3).Value = 0.306821227259698
$ws.Cells.Item(33, 4).Value = 0.7527432677738641
$ws.Cells.Item(33, 5).Value = 10.19245300693656
$ws.Cells.Item(33, 7).Value = 11.91266994300607
$ws.Cells.Item(34, 2).Value = 3.286832544864788
$ws.Cells.Item(34, 3).Value = 1.655778082260271
$ws.Cells.Item(34, 4).Value = 0.7527432677738641
$ws.Cells.Item(34, 5).Value = 0.4942365360607697
$ws.Cells.Item(34, 7).Value = 6.189590430959694
$ws.Cells.Item(35, 2).Value = 3.286832544864788
$ws.Cells.Item(35, 3).Value = 1.655778082260271
$ws.Cells.Item(35, 4).Value = 0.7527432677738641
$ws.Cells.Item(35, 5).Value = 0.4942365360607697
$ws.Cells.Item(35, 7).Value = 6.189590430959694
$ws.Cells.Item(36, 2).Value = 3.286832544864788
$ws.Cells.Item(36, 3).Value = 1.655778082260271
$ws.Cells.Item(36, 4).Value = 0.7527432677738641
$ws.Cells.Item(36, 5).Value = 0.4942365360607697
$ws.Cells.Item(36, 7).Value = 6.189590430959694
$ws.Cells.Item(37, 2).Value = 3.286832544864788
$ws.Cells.Item(37, 3).Value = 1.655778082260271
$ws.Cells.Item(37, 4).Value = 0.7527432677738641
$ws.Cells.Item(37, 5).Value = 0.4942365360607697
$ws.Cells.Item(37, 7).Value = 6.189590430959694
$ws.Cells.Item(38, 2).Value = 3.286832544864788
$ws.Cells.Item(38, 3).Value = 1.655778082260271
$ws.Cells.Item(38, 4).Value = 261.3203778131603
$ws.Cells.Item(38, 5).Value = 10.19245300693656
$ws.Cells.Item(38, 7).Value = 276.4554414472219
$ws.Cells.Item(39, 2).Value = 3.286832544864788
$ws.Cells.Item(39, 3).Value = 1.655778082260271
$ws.Cells.Item(39, 4).Value = 0.7527432677738641
$ws.Cells.Item(39, 5).Value = 0.4942365360607697
$ws.Cells.Item(39, 7).Value = 6.189590430959694
$ws.Cells.Item(40, 2).Value = 3.286832544864788
$ws.Cells.Item(40, 3).Value = 1.655778082260271
$ws.Cells.Item(40, 4).Value = 0.7527432677738641
$ws.Cells.Item(40, 5).Value = 0.4942365360607697
$ws.Cells.Item(40, 7).Value = 6.189590430959694
